$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$numericCells = @('D5','D6','D9','D11','D15','D19','D20','D21','D22','D23','D24','D28','D29','D33','D34','D36','D38','D41','D43','D44','D46','D48','D51')
foreach ($addr in $numericCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range('D2').Value = '66.522.86'
$ws.Range('E2').Value = '  -4.84%  '
$ws.Range('D3').Value = '3.444.90'
$ws.Range('E3').Value = '  -6.65%  '
$ws.Range('D5').Value = '598.86'
$ws.Range('E5').Value = '  -7.61%  '
$ws.Range('D6').Value = '146.58'
$ws.Range('E6').Value = '  -9.54%  '
$ws.Range('D7').Value = '3.444.47'
$ws.Range('E7').Value = '  -6.69%  '
$ws.Range('D9').Value = '0.475'
$ws.Range('E9').Value = '  -5.27%  '
$ws.Range('E10').Value = '  -7.50%  '
$ws.Range('D11').Value = '6.85'
$ws.Range('E11').Value = '  -4.95%  '
$ws.Range('E12').Value = '  -6.37%  '
$ws.Range('E13').Value = '  -9.08%  '
$ws.Range('D14').Value = '4.027.35'
$ws.Range('E14').Value = '  -6.59%  '
$ws.Range('D15').Value = '30.76'
$ws.Range('E15').Value = '  -6.12%  '
$ws.Range('B16').Value = 'WrappedBTC'
$ws.Range('C16').Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range('D16').Value = '66.458.40'
$ws.Range('E16').Value = '  -4.87%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.435.76'
$ws.Range('E17').Value = '  -6.36%  '
$ws.Range('E18').Value = '  -1.03%  '
$ws.Range('D19').Value = '6.33'
$ws.Range('E19').Value = '  -2.99%  '
$ws.Range('D20').Value = '14.76'
$ws.Range('E20').Value = '  -7.74%  '
$ws.Range('D21').Value = '436.92'
$ws.Range('E21').Value = '  -7.35%  '
$ws.Range('D22').Value = '8.90'
$ws.Range('E22').Value = '  -14.41%  '
$ws.Range('D23').Value = '0.616'
$ws.Range('D24').Value = '76.59'
$ws.Range('E24').Value = '  -4.31%  '
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('D26').Value = '3.582.83'
$ws.Range('E27').Value = '  -4.28%  '
$ws.Range('D28').Value = '9.92'
$ws.Range('E28').Value = '  -9.52%  '
$ws.Range('D29').Value = '8.17'
$ws.Range('E29').Value = '  -10.57%  '
$ws.Range('E30').Value = '  -6.19%  '
$ws.Range('E31').Value = '  -0.01%  '
$ws.Range('E32').Value = '  -11.17%  '
$ws.Range('D33').Value = '0.159'
$ws.Range('E33').Value = '  -5.79%  '
$ws.Range('D34').Value = '25.33'
$ws.Range('E34').Value = '  -5.38%  '
$ws.Range('E35').Value = '  -7.64%  '
$ws.Range('D36').Value = '1.82'
$ws.Range('E36').Value = '  -9.21%  '
$ws.Range('D37').Value = '3.429.83'
$ws.Range('E37').Value = '  -6.94%  '
$ws.Range('D38').Value = '7.81'
$ws.Range('E38').Value = '  -7.57%  '
$ws.Range('E39').Value = '  +0.06%  '
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('D41').Value = '173.18'
$ws.Range('E41').Value = '  -3.80%  '
$ws.Range('E42').Value = '  -5.58%  '
$ws.Range('D43').Value = '0.0853'
$ws.Range('E43').Value = '  -5.63%  '
$ws.Range('D44').Value = '5.35'
$ws.Range('E44').Value = '  -9.24%  '
$ws.Range('E45').Value = '  -6.67%  '
$ws.Range('D46').Value = '45.22'
$ws.Range('E46').Value = '  -3.01%  '
$ws.Range('E47').Value = '  -4.33%  '
$ws.Range('D48').Value = '25.77'
$ws.Range('E48').Value = '  -12.30%  '
$ws.Range('E49').Value = '  -5.05%  '
$ws.Range('E50').Value = '  -14.56%  '
$ws.Range('D51').Value = '0.988'
$ws.Range('E51').Value = '  -6.56%  '

foreach ($addr in $numericCells) {
    $ws.Range($addr).Style = "Normal"
}
